# Update "想去人数" (people interested) counts on both the "展览" sheet
# and the "全部类型" sheet (which mirrors the same rows).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 43
    $ws.Range("F3").Value = 143
    $ws.Range("F4").Value = 16
    $ws.Range("F5").Value = 43
}
